# Added 1.1.0 of term
# Bumps the ValueSet "Version" metadata cell and its "Date" cell on the
# Metadata sheet to reflect the new release.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Property/Value table: row 3 = Version, row 8 = Date
$ws.Range("B3").Value = "1.1.0"
$ws.Range("B8").Value = "2023-07-10T23:08:03+02:00"
